$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("B1").Value = 0.92313574545218779
$ws.Range("BP1").Value = 0.55788142743996483
$ws.Range("D2").Value = 0.97910195997093075
$ws.Range("BL2").Value = 0.93680640711625585
$ws.Range("B3").Value = 0.83856752585377681
$ws.Range("C4").Value = 0.88551897155931614
$ws.Range("E4").Value = 0.78498233326063604
$ws.Range("C5").Value = 0.99624075975271387
$ws.Range("F5").Value = 0.86094768457230919
$ws.Range("G5").Value = 0.99341872582772472
$ws.Range("D6").Value = 0.89847954889289261
$ws.Range("F7").Value = 0.89093090566269906
$ws.Range("H7").Value = 0.94929052201329844
$ws.Range("I7").Value = 0.96783627251711946
$ws.Range("AJ8").Value = 0.8012167648932973
$ws.Range("H9").Value = 0.88883777320996027
$ws.Range("K9").Value = 0.9329625516051685
$ws.Range("H10").Value = 0.82695034705679504
$ws.Range("K10").Value = 0.90009525489602549
$ws.Range("L11").Value = 0.61093629266642036
$ws.Range("J12").Value = 0.91702987015885074
$ws.Range("N12").Value = 0.86687058118945537
$ws.Range("BK12").Value = 0.86525348616341824
$ws.Range("K13").Value = 0.93067633083854018
$ws.Range("N13").Value = 0.80895346825588987
$ws.Range("O13").Value = 0.74762509356371876
$ws.Range("BM13").Value = 0.75495300125139542
$ws.Range("O14").Value = 0.98040116093289698
$ws.Range("O16").Value = 0.99372039407045554
$ws.Range("R16").Value = 0.705872864617199
$ws.Range("M17").Value = 0.71464953594163827
$ws.Range("O17").Value = 0.86700726123233651
$ws.Range("P17").Value = 0.76595188947004478
$ws.Range("Q18").Value = 0.63350814238406028
$ws.Range("S18").Value = 0.92063201738701461
$ws.Range("T18").Value = 0.79398392157762188
$ws.Range("U19").Value = 0.75973351053034288
$ws.Range("AS19").Value = 0.70907959123338959
$ws.Range("S20").Value = 0.96914794883188948
$ws.Range("V20").Value = 0.87417628998175045
$ws.Range("AB20").Value = 0.75119602548630249
$ws.Range("V21").Value = 0.64410176689545695
$ws.Range("A22").Value = 0.99737648466367634
$ws.Range("W22").Value = 0.73743284134470755
$ws.Range("X22").Value = 0.66498218831386757
$ws.Range("U23").Value = 0.9144470888240569
$ws.Range("X23").Value = 0.96086550672178994
$ws.Range("Y24").Value = 0.60931809439180684
$ws.Range("Z24").Value = 0.94194762869137616
$ws.Range("AW24").Value = 0.86316096240594364
$ws.Range("Z25").Value = 0.81898223032978623
$ws.Range("AA25").Value = 0.77482042884501467
$ws.Range("AB26").Value = 0.93452575603072696
$ws.Range("BB26").Value = 0.99819040904943868
$ws.Range("AB27").Value = 0.59600798492247065
$ws.Range("AC27").Value = 0.86138855687618854
$ws.Range("AT27").Value = 0.94120745110295023
$ws.Range("AD28").Value = 0.85386237916951147
$ws.Range("AD29").Value = 0.97024706400028138
$ws.Range("AT29").Value = 0.91907819988094541
$ws.Range("AF30").Value = 0.85418539133201488
$ws.Range("AD31").Value = 0.90727201843792349
$ws.Range("AG31").Value = 0.91710928211096743
$ws.Range("AE32").Value = 0.98976673468713305
$ws.Range("AP32").Value = 0.79421502247427922
$ws.Range("BF32").Value = 0.81131106235094341
$ws.Range("AF34").Value = 0.91528435633106864
$ws.Range("AG34").Value = 0.81491194336436368
$ws.Range("AG35").Value = 0.64195095347579501
$ws.Range("AH35").Value = 0.9964422896014753
$ws.Range("AH36").Value = 0.72690288460775543
$ws.Range("AI36").Value = 0.91359205665089316
$ws.Range("W37").Value = 0.8702033283977193
$ws.Range("AI37").Value = 0.76251957947865612
$ws.Range("AJ37").Value = 0.54015268984339782
$ws.Range("AK38").Value = 0.88755253083480889
$ws.Range("AK39").Value = 0.9422428733060394
$ws.Range("AL39").Value = 0.8884667776311137
$ws.Range("AN39").Value = 0.69475343007811219
$ws.Range("AP40").Value = 0.75804059274254554
$ws.Range("AM41").Value = 0.79287572136826268
$ws.Range("AN41").Value = 0.65074387770051878
$ws.Range("BK41").Value = 0.89652808142175533
$ws.Range("AQ42").Value = 0.74150216665540603
$ws.Range("AR42").Value = 0.88795322596437587
$ws.Range("AR43").Value = 0.92972324492604352
$ws.Range("L44").Value = 0.81863693612550981
$ws.Range("BD44").Value = 0.67313761501525704
$ws.Range("AQ45").Value = 0.66491977194609886
$ws.Range("AR45").Value = 0.74154718462812785
$ws.Range("AT45").Value = 0.89330025651321532
$ws.Range("F46").Value = 0.94223458718414432
$ws.Range("AL46").Value = 0.60347639580491663
$ws.Range("AR46").Value = 0.85044747924654474
$ws.Range("AT47").Value = 0.97019962169291141
$ws.Range("AW47").Value = 0.86976496778881551
$ws.Range("AU48").Value = 0.95316189715857869
$ws.Range("AV49").Value = 0.66491782777649078
$ws.Range("AX49").Value = 0.67068220069691353
$ws.Range("AV50").Value = 0.825771065131916
$ws.Range("AZ50").Value = 0.81079347382780043
$ws.Range("AX51").Value = 0.74851613572663034
$ws.Range("BB51").Value = 0.84424014695176419
$ws.Range("M52").Value = 0.82748492132668761
$ws.Range("BB53").Value = 0.63317288750863643
$ws.Range("AH54").Value = 0.66162696887393069
$ws.Range("AZ54").Value = 0.95643153401352521
$ws.Range("BC54").Value = 0.91597290275283005
$ws.Range("N55").Value = 0.9867010135224058
$ws.Range("BA55").Value = 0.92251887362410678
$ws.Range("BD55").Value = 0.94764378852858899
$ws.Range("T57").Value = 0.83380321146399128
$ws.Range("Y57").Value = 0.86154274821745602
$ws.Range("BD57").Value = 0.76508513596687089
$ws.Range("BG57").Value = 0.9493881740094664
$ws.Range("I58").Value = 0.86132080176478032
$ws.Range("BG58").Value = 0.85047636527644399
$ws.Range("BH59").Value = 0.68879096978397691
$ws.Range("BI59").Value = 0.63533892607262565
$ws.Range("BI60").Value = 0.81477944372211653
$ws.Range("BJ61").Value = 0.6951669892977439
$ws.Range("BK61").Value = 0.87337313634082503
$ws.Range("BH62").Value = 0.8849259433337815
$ws.Range("BL62").Value = 0.7845480829184629
$ws.Range("N63").Value = 0.8158625354688458
$ws.Range("AY63").Value = 0.98580328887000346
$ws.Range("BJ63").Value = 0.96249159245417304
$ws.Range("BK64").Value = 0.97073660248045823
$ws.Range("BN65").Value = 0.74407013836354841
$ws.Range("AO66").Value = 0.80447689598960337
$ws.Range("BL66").Value = 0.62369917096386507
$ws.Range("BO66").Value = 0.66944139417446324
$ws.Range("BP66").Value = 0.83816465174518129
$ws.Range("A67").Value = 0.75225387665590859
$ws.Range("U68").Value = 0.88156477644096065
$ws.Range("BO68").Value = 0.90802041004423706
